# Applies the scheduled-runner value refresh to the Louisoix_Profits workbook.
# Each touched row keeps its Leve metadata (A-G) intact; only the live Market
# Board snapshot columns (H-N: currentAveragePrice.. LeveProfitHQ) are updated
# to the latest pull. A few rows gain/lose a trailing HQ-profit cell (N, or M)
# depending on whether an HQ price now exists for that item.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 289.6
$ws.Range("I12").Value = 237.25
$ws.Range("K12").Value = 237.25
$ws.Range("M12").Value = -67.25
$ws.Range("H18").Value = 449.5
$ws.Range("I18").Value = 449.5
$ws.Range("K18").Value = 449.5
$ws.Range("M18").Value = -165.5
$ws.Range("H33").Value = 8812.77
$ws.Range("I33").Value = 8547.166999999999
$ws.Range("J33").Value = 12000
$ws.Range("K33").Value = 8547.166999999999
$ws.Range("L33").Value = 12000
$ws.Range("M33").Value = -8318.166999999999
$ws.Range("N33").Value = -12458
$ws.Range("H40").Value = 3371.2856
$ws.Range("I40").Value = 1716.6666
$ws.Range("J40").Value = 4612.25
$ws.Range("K40").Value = 1716.6666
$ws.Range("L40").Value = 4612.25
$ws.Range("M40").Value = -1541.6666
$ws.Range("N40").Value = -4962.25
$ws.Range("H55").Value = 595.09375
$ws.Range("I55").Value = 412.75
$ws.Range("K55").Value = 412.75
$ws.Range("M55").Value = -198.75
$ws.Range("H88").Value = 1437.5714
$ws.Range("J88").Value = 1338.8
$ws.Range("L88").Value = 1338.8
$ws.Range("N88").Value = -2150.8
$ws.Range("H91").Value = 1437.5714
$ws.Range("J91").Value = 1338.8
$ws.Range("L91").Value = 1338.8
$ws.Range("N91").Value = -4146.8
$ws.Range("H100").Value = 2654.5
$ws.Range("I100").Value = 1367.0834
$ws.Range("J100").Value = 3941.9167
$ws.Range("K100").Value = 1367.0834
$ws.Range("L100").Value = 3941.9167
$ws.Range("M100").Value = -826.0834
$ws.Range("N100").Value = -5023.9167
$ws.Range("H111").Value = 2224.5557
$ws.Range("I111").Value = 2654.3333
$ws.Range("J111").Value = 1365
$ws.Range("K111").Value = 7962.999899999999
$ws.Range("L111").Value = 4095
$ws.Range("M111").Value = -4895.999899999999
$ws.Range("N111").Value = -10229
$ws.Range("H113").Value = 11680.857
$ws.Range("I113").Value = 19990
$ws.Range("K113").Value = 19990
$ws.Range("M113").Value = -16736
$ws.Range("H115").Value = 1342.5385
$ws.Range("I115").Value = 1099.6
$ws.Range("K115").Value = 3298.8
$ws.Range("M115").Value = -1731.8
$ws.Range("H125").Value = 1326.375
$ws.Range("I125").Value = 925
$ws.Range("K125").Value = 8325
$ws.Range("M125").Value = -5865
$ws.Range("H127").Value = 1430.5
$ws.Range("I127").Value = 1635.6666
$ws.Range("K127").Value = 4906.9998
$ws.Range("M127").Value = 53.0002000000004
$ws.Range("H135").Value = 1587.2727
$ws.Range("J135").Value = 966
$ws.Range("L135").Value = 8694
$ws.Range("N135").Value = -13764
$ws.Range("H137").Value = 1845.9166
$ws.Range("I137").Value = 1044.875
$ws.Range("K137").Value = 3134.625
$ws.Range("M137").Value = -584.625
$ws.Range("H138").Value = 2442.375
$ws.Range("I138").Value = 2963.5
$ws.Range("J138").Value = 2070.1428
$ws.Range("K138").Value = 8890.5
$ws.Range("L138").Value = 6210.428400000001
$ws.Range("M138").Value = -3750.5
$ws.Range("N138").Value = -16490.4284
$ws.Range("H141").Value = 7040.857
$ws.Range("I141").Value = 7040.857
$ws.Range("K141").Value = 21122.571
$ws.Range("M141").Value = -15942.571

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29715.928
$ws.Range("I32").Value = 34382.484
$ws.Range("J32").Value = 10466.375
$ws.Range("K32").Value = 34382.484
$ws.Range("L32").Value = 10466.375
$ws.Range("M32").Value = -34095.484
$ws.Range("N32").Value = -11040.375
$ws.Range("H45").Value = 4212.522
$ws.Range("I45").Value = 3296.2856
$ws.Range("J45").Value = 5637.778
$ws.Range("K45").Value = 3296.2856
$ws.Range("L45").Value = 5637.778
$ws.Range("M45").Value = -2919.2856
$ws.Range("N45").Value = -6391.778
$ws.Range("H88").Value = 1811.8235
$ws.Range("I88").Value = 1124.75
$ws.Range("J88").Value = 2023.2307
$ws.Range("K88").Value = 1124.75
$ws.Range("L88").Value = 2023.2307
$ws.Range("M88").Value = -718.75
$ws.Range("N88").Value = -2835.2307
$ws.Range("H91").Value = 1811.8235
$ws.Range("I91").Value = 1124.75
$ws.Range("J91").Value = 2023.2307
$ws.Range("K91").Value = 1124.75
$ws.Range("L91").Value = 2023.2307
$ws.Range("M91").Value = 279.25
$ws.Range("N91").Value = -4831.2307
$ws.Range("H97").Value = 13920.9
$ws.Range("I97").Value = 21582
$ws.Range("K97").Value = 21582
$ws.Range("M97").Value = -21086

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1901.4138
$ws.Range("I20").Value = 1734.7
$ws.Range("K20").Value = 1734.7
$ws.Range("M20").Value = -1487.7
$ws.Range("H102").Value = 10879
$ws.Range("I102").Value = 10879
$ws.Range("K102").Value = 10879
$ws.Range("M102").Value = -7634
$ws.Range("H107").Value = 1038.7333
$ws.Range("I107").Value = 755.7857
$ws.Range("K107").Value = 755.7857
$ws.Range("M107").Value = 1164.2143
$ws.Range("H134").Value = 2235.8667
$ws.Range("I134").Value = 1700.8334
$ws.Range("J134").Value = 4376
$ws.Range("K134").Value = 5102.5002
$ws.Range("L134").Value = 13128
$ws.Range("M134").Value = -2567.5002
$ws.Range("N134").Value = -18198

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 44552
$ws.Range("J43").Value = 44552
$ws.Range("L43").Value = 44552
$ws.Range("N43").Value = -44920
$ws.Range("H88").Value = 17188.715
$ws.Range("I88").Value = 4995
$ws.Range("J88").Value = 19221
$ws.Range("K88").Value = 4995
$ws.Range("L88").Value = 19221
$ws.Range("M88").Value = -4589
$ws.Range("N88").Value = -20033
$ws.Range("H91").Value = 17188.715
$ws.Range("I91").Value = 4995
$ws.Range("J91").Value = 19221
$ws.Range("K91").Value = 4995
$ws.Range("L91").Value = 19221
$ws.Range("M91").Value = -3591
$ws.Range("N91").Value = -22029
$ws.Range("H101").Value = 44552
$ws.Range("J101").Value = 44552
$ws.Range("L101").Value = 44552
$ws.Range("N101").Value = -51042

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1123.875
$ws.Range("I8").Value = 1123.875
$ws.Range("K8").Value = 3371.625
$ws.Range("M8").Value = -3232.625
$ws.Range("H128").Value = 135070
$ws.Range("I128").Value = 135070
$ws.Range("K128").Value = 405210
$ws.Range("M128").Value = -400230
$ws.Range("H129").Value = 718582.8
$ws.Range("J129").Value = 2834966.2
$ws.Range("L129").Value = 8504898.600000001
$ws.Range("N129").Value = -8514898.600000001
$ws.Range("H131").Value = 3852302
$ws.Range("I131").Value = 18700.5
$ws.Range("J131").Value = 5002382.5
$ws.Range("K131").Value = 56101.5
$ws.Range("L131").Value = 15007147.5
$ws.Range("M131").Value = -51061.5
$ws.Range("N131").Value = -15017227.5
$ws.Range("H139").Value = 1378.1428
$ws.Range("I139").Value = 1378.1428
$ws.Range("K139").Value = 4134.428400000001
$ws.Range("M139").Value = 1005.571599999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2752.75
$ws.Range("I102").Value = 1670.6666
$ws.Range("J102").Value = 5999
$ws.Range("K102").Value = 1670.6666
$ws.Range("L102").Value = 5999
$ws.Range("M102").Value = -48.66660000000002
$ws.Range("N102").Value = -9243
$ws.Range("H122").Value = 2056.1667
$ws.Range("I122").Value = 1635.7142
$ws.Range("J122").Value = 4999.3335
$ws.Range("K122").Value = 4907.142599999999
$ws.Range("L122").Value = 14998.0005
$ws.Range("M122").Value = -2457.142599999999
$ws.Range("N122").Value = -19898.0005
$ws.Range("H126").Value = 7319.25
$ws.Range("I126").Value = 6892.25
$ws.Range("J126").Value = 7746.25
$ws.Range("K126").Value = 20676.75
$ws.Range("L126").Value = 23238.75
$ws.Range("M126").Value = -18206.75
$ws.Range("N126").Value = -28178.75
$ws.Range("H132").Value = 169491.5
$ws.Range("I132").Value = 252537.25
$ws.Range("J132").Value = 3400
$ws.Range("K132").Value = 757611.75
$ws.Range("L132").Value = 10200
$ws.Range("M132").Value = -755081.75
$ws.Range("N132").Value = -15260

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 13577.2
$ws.Range("I46").Value = 19810.234
$ws.Range("J46").Value = 5426.3076
$ws.Range("K46").Value = 19810.234
$ws.Range("L46").Value = 5426.3076
$ws.Range("M46").Value = -19622.234
$ws.Range("N46").Value = -5802.3076

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 24039
$ws.Range("I49").Value = 24039
$ws.Range("K49").Value = 24039
$ws.Range("M49").Value = -23809
$ws.Range("H81").Value = 2173.4285
$ws.Range("I81").Value = 1778.5
$ws.Range("J81").Value = 2700
$ws.Range("K81").Value = 3557
$ws.Range("L81").Value = 5400
$ws.Range("M81").Value = -2496
$ws.Range("N81").Value = -7522
$ws.Range("H84").Value = 2173.4285
$ws.Range("I84").Value = 1778.5
$ws.Range("J84").Value = 2700
$ws.Range("K84").Value = 17785
$ws.Range("L84").Value = 27000
$ws.Range("M84").Value = -12481
$ws.Range("N84").Value = -37608
$ws.Range("H107").Value = 1045.75
$ws.Range("I107").Value = 1045.75
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3137.25
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1217.25
$ws.Range("N107").ClearContents()
$ws.Range("H132").Value = 108999.5
$ws.Range("J132").Value = 1178
$ws.Range("L132").Value = 3534
$ws.Range("N132").Value = -8594

Write-Output "Updated 242 cell(s); cleared 1 cell(s)."
